# "update phan cong 2"
# The completion % for the 3rd task (row 4) dropped from 100% (1) to 80% (0.8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 0.8

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("C9").Select() | Out-Null
